# Revise average and remove unnecessary average in Summary124
# Updates the R:W columns for rows 11, 12, 13, 15, 16 with revised values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("12.4")

# Row 11
$ws.Range("R11").Value = 0.030974471933338
$ws.Range("S11").Value = 0.0068479122191713
$ws.Range("T11").Value = 0.052953164575325
$ws.Range("U11").Value = 0.0016457306606764
$ws.Range("V11").Value = 0.041963818254331
$ws.Range("W11").Value = 0.0042468214399239

# Row 12
$ws.Range("R12").Value = 0.73876409792649
$ws.Range("S12").Value = 0.30994441750723
$ws.Range("T12").Value = 2.0428731993907
$ws.Range("U12").Value = 0.064585360230034
$ws.Range("V12").Value = 1.3908186486586
$ws.Range("W12").Value = 0.18726488886863

# Row 13
$ws.Range("R13").Value = 0.027401020605904
$ws.Range("S13").Value = 0.0065363390157085
$ws.Range("T13").Value = 0.057164621745122
$ws.Range("U13").Value = 0.0017988218849254
$ws.Range("V13").Value = 0.042282821175513
$ws.Range("W13").Value = 0.004167580450317

# Row 15
$ws.Range("R15").Value = 4.3135256908468
$ws.Range("S15").Value = 0.68471867233326
$ws.Range("T15").Value = 9.985659882498
$ws.Range("U15").Value = 0.31718588024084
$ws.Range("V15").Value = 7.1495927866724
$ws.Range("W15").Value = 0.50095227628705

# Row 16
$ws.Range("R16").Value = 1.3355304696079
$ws.Range("S16").Value = 0.34583884439391
$ws.Range("T16").Value = 4.0207474287372
$ws.Range("U16").Value = 0.143824413826
$ws.Range("V16").Value = 2.6781389491726
$ws.Range("W16").Value = 0.24483162910996
